# Refresh cryptocurrency Price (D) and Volume(1h) (E) figures for sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column values that look like plain numbers need a leading apostrophe so
# Excel keeps storing them as text (matching the source data, which mixes
# European-style grouped numbers with plain decimals in the same column).

$ws.Range("D2").Value = '64.300.48'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '3.504.04'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''586.22'
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("D6").Value = '''134.39'
$ws.Range("E6").Value = '  +3.74%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  +1.17%  '
$ws.Range("E9").Value = '  +2.00%  '
$ws.Range("E10").Value = '  +2.47%  '
$ws.Range("E11").Value = '  +3.09%  '
$ws.Range("D12").Value = '4.099.02'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").Value = '''0.0000182'
$ws.Range("E13").Value = '  +4.58%  '
$ws.Range("E14").Value = '  +1.42%  '
$ws.Range("D15").Value = '3.505.00'
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").Value = '''26.00'
$ws.Range("E16").Value = '  -3.75%  '
$ws.Range("D17").Value = '64.310.27'
$ws.Range("E17").Value = '  +1.29%  '
$ws.Range("D18").Value = '''9.91'
$ws.Range("E18").Value = '  +1.43%  '
$ws.Range("E19").Value = '  +3.36%  '
$ws.Range("D20").Value = '''13.67'
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("D21").Value = '''393.64'
$ws.Range("E21").Value = '  +3.96%  '
$ws.Range("E22").Value = '  +0.57%  '
$ws.Range("D23").Value = '3.643.40'
$ws.Range("E23").Value = '  +0.78%  '
$ws.Range("D24").Value = '''74.28'
$ws.Range("E24").Value = '  +1.74%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  +2.38%  '
$ws.Range("E27").Value = '  -0.19%  '
$ws.Range("D28").Value = '''7.44'
$ws.Range("E28").Value = '  +0.64%  '
$ws.Range("D29").Value = '''1.50'
$ws.Range("E29").Value = '  -3.29%  '
$ws.Range("D30").Value = '''8.30'
$ws.Range("E30").Value = '  +1.93%  '
$ws.Range("D31").Value = '''2.24'
$ws.Range("E31").Value = '  +1.50%  '
$ws.Range("D32").Value = '3.523.15'
$ws.Range("E32").Value = '  +1.10%  '
$ws.Range("E33").Value = '  +5.01%  '
$ws.Range("D35").Value = '''23.47'
$ws.Range("E35").Value = '  +1.02%  '
$ws.Range("E36").Value = '  -0.79%  '
$ws.Range("D37").Value = '''1.57'
$ws.Range("E37").Value = '  +1.88%  '
$ws.Range("E38").Value = '  +0.82%  '
$ws.Range("D39").Value = '''163.86'
$ws.Range("E39").Value = '  +2.57%  '
$ws.Range("D40").Value = '''0.0785'
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").Value = '''25.13'
$ws.Range("E43").Value = '  -2.95%  '
$ws.Range("E44").Value = '  +2.45%  '
$ws.Range("E45").Value = '  +3.70%  '
$ws.Range("E46").Value = '  -1.67%  '
$ws.Range("D47").Value = '2.461.97'
$ws.Range("E47").Value = '  +2.12%  '
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("D49").Value = '''0.900'
$ws.Range("E49").Value = '  +2.12%  '
$ws.Range("D50").Value = '''0.0262'
$ws.Range("E50").Value = '  -0.21%  '
$ws.Range("E51").Value = '  +0.52%  '
